# 2023-08-23 update next cosmetic
#
# The backlog's oldest/next-up shop item (id 507, "Nyakuza Mask" by
# snekiecr8) has been resolved and is removed from the sheet; the sheet
# now reflects the next item in line (id 1027, "Count Olaf" cosmetic).
#
# Structurally this:
#   - drops the old helper column A (the previous "current item id" cell
#     that all the other columns were offset around), shifting B:M left
#     into A:L
#   - removes the last row (the old id-507/Nyakuza Mask/snekiecr8 lookup
#     row), shifting the used range from 17 rows to 16
#   - makes room for two more trailing helper columns (M:N) by inserting
#     two blank columns, which also pushes the little
#     "shop item rows / citb user(s) / citb comment" lookup table from
#     J:L out to L:N
#   - relabels the header row back to the generic col_0..col_13 scheme
#   - fills in the newly exposed helper cells

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stale "current item" helper column; everything shifts one
# column to the left (old B -> A, old C -> B, ... old M -> L).
$ws.Columns.Item(1).Delete()

# Drop the stale id-507 lookup row at the bottom.
$ws.Rows.Item(17).Delete()

# Make room for two more helper columns at the right of the small
# "shop item rows / citb user(s) / citb comment" lookup table, pushing
# it from J:L out to L:N.
$ws.Columns.Item(10).Insert()
$ws.Columns.Item(10).Insert()

# Header row: generic col_N labels.
$ws.Range("A1").Value = "row_number"
$ws.Range("B1").Value = "col_1"
$ws.Range("C1").Value = "Type to search:"
$ws.Range("D1").Value = "col_3"
$ws.Range("E1").Value = "col_4"
$ws.Range("F1").Value = "col_5"
$ws.Range("G1").Value = "col_6"
$ws.Range("H1").Value = "col_7"
$ws.Range("I1").Value = "col_8"
$ws.Range("J1").Value = "col_9"
$ws.Range("K1").Value = "col_10"
$ws.Range("L1").Value = "col_11"
$ws.Range("M1").Value = "col_12"
$ws.Range("N1").Value = "col_13"

# The two newly-inserted helper columns are blank on every data row.
$ws.Range("J2:K16").Value = ""

# "Shop candidates" selector now points at row 1 (was 2).
$ws.Range("D4").Value = 1
